# Scheduled market-data refresh: update cached price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the Leve sheets.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 413.58334
$ws.Range("I28").Value = 96.59999999999999
$ws.Range("K28").Value = 96.59999999999999
$ws.Range("M28").Value = 388.4
# Row 41 (Leve Item ID 5478)
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()
# Row 42 (Leve Item ID 4600)
$ws.Range("H42").Value = 330.45456
$ws.Range("I42").Value = 204.71428
$ws.Range("K42").Value = 614.14284
$ws.Range("M42").Value = -384.14284
# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 6750
$ws.Range("I43").Value = 5500
$ws.Range("J43").Value = 8000
$ws.Range("K43").Value = 5500
$ws.Range("L43").Value = 8000
$ws.Range("M43").Value = -5431
$ws.Range("N43").Value = -8138
# Row 55 (Leve Item ID 5517)
$ws.Range("H55").Value = 1015.5714
$ws.Range("I55").Value = 1513.2
$ws.Range("J55").Value = 563.1818
$ws.Range("K55").Value = 1513.2
$ws.Range("L55").Value = 563.1818
$ws.Range("M55").Value = -1299.2
$ws.Range("N55").Value = -991.1818
# Row 118 (Leve Item ID 27958)
$ws.Range("H118").Value = 160.66667
$ws.Range("I118").Value = 152.8
$ws.Range("J118").Value = 200
$ws.Range("K118").Value = 458.4
$ws.Range("L118").Value = 600
$ws.Range("M118").Value = 1198.6
$ws.Range("N118").Value = -3914

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5 (Leve Item ID 5091)
$ws.Range("H5").Value = 195.13333
$ws.Range("I5").Value = 204.14285
$ws.Range("J5").Value = 69
$ws.Range("K5").Value = 204.14285
$ws.Range("L5").Value = 69
$ws.Range("M5").Value = -92.14285000000001
$ws.Range("N5").Value = -293
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 5572.6
$ws.Range("I32").Value = 5572.6
$ws.Range("K32").Value = 5572.6
$ws.Range("M32").Value = -5285.6

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4 (Leve Item ID 5091)
$ws.Range("H4").Value = 195.13333
$ws.Range("I4").Value = 204.14285
$ws.Range("J4").Value = 69
$ws.Range("K4").Value = 204.14285
$ws.Range("L4").Value = 69
$ws.Range("M4").Value = -89.14285000000001
$ws.Range("N4").Value = -299
# Row 141 (Leve Item ID 43278)
$ws.Range("H141").Value = 60000
$ws.Range("I141").Value = 20000
$ws.Range("K141").Value = 20000
$ws.Range("M141").Value = -14820

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 1173.6154
$ws.Range("I16").Value = 1047.2
$ws.Range("K16").Value = 1047.2
$ws.Range("M16").Value = -760.2
# Row 22 (Leve Item ID 5367)
$ws.Range("H22").Value = 992.2857
$ws.Range("I22").Value = 957.6667
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 957.6667
$ws.Range("L22").Value = 1200
$ws.Range("M22").Value = -607.6667
$ws.Range("N22").Value = -1900
# Row 50 (Leve Item ID 1862)
$ws.Range("H50").Value = 20000
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 2502500
$ws.Range("I99").Value = 2502500
$ws.Range("K99").Value = 2502500
$ws.Range("M99").Value = -2501002
# Row 100 (Leve Item ID 34388)
$ws.Range("H100").Value = 100780
$ws.Range("J100").Value = 100780
$ws.Range("L100").Value = 100780
$ws.Range("N100").Value = -102944
# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 583.1111
$ws.Range("I107").Value = 624.5
$ws.Range("J107").Value = 571.2857
$ws.Range("K107").Value = 624.5
$ws.Range("L107").Value = 571.2857
$ws.Range("M107").Value = 1295.5
$ws.Range("N107").Value = -4411.2857
# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 1173.6154
$ws.Range("I113").Value = 1047.2
$ws.Range("K113").Value = 1047.2
$ws.Range("M113").Value = 1122.8
# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 2502500
$ws.Range("I126").Value = 2502500
$ws.Range("K126").Value = 7507500
$ws.Range("M126").Value = -7505030

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 23 (Leve Item ID 4858)
$ws.Range("H23").Value = 107.68421
$ws.Range("I23").Value = 38.8
$ws.Range("J23").Value = 132.28572
$ws.Range("K23").Value = 116.4
$ws.Range("L23").Value = 396.85716
$ws.Range("M23").Value = 118.6
$ws.Range("N23").Value = -866.85716
# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 1868.909
$ws.Range("J68").Value = 1973
$ws.Range("L68").Value = 5919
$ws.Range("N68").Value = -7541
# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 1868.909
$ws.Range("J71").Value = 1973
$ws.Range("L71").Value = 17757
$ws.Range("N71").Value = -25869
# Row 80 (Leve Item ID 12890)
$ws.Range("H80").Value = 5000
$ws.Range("I80").Value = 5000
$ws.Range("K80").Value = 15000
$ws.Range("M80").Value = -14064
# Row 83 (Leve Item ID 12890)
$ws.Range("H83").Value = 5000
$ws.Range("I83").Value = 5000
$ws.Range("K83").Value = 45000
$ws.Range("M83").Value = -40320
# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 767.2
$ws.Range("I132").Value = 489.6
$ws.Range("J132").Value = 1044.8
$ws.Range("K132").Value = 4406.400000000001
$ws.Range("L132").Value = 9403.199999999999
$ws.Range("M132").Value = -1876.400000000001
$ws.Range("N132").Value = -14463.2
# Row 137 (Leve Item ID 44088)
$ws.Range("H137").Value = 1000
$ws.Range("I137").Value = 1000
$ws.Range("K137").Value = 3000
$ws.Range("M137").Value = 2100
# Row 141 (Leve Item ID 44076)
$ws.Range("H141").Value = 1286.8889
$ws.Range("I141").Value = 1098.1428
$ws.Range("K141").Value = 3294.4284
$ws.Range("M141").Value = 1885.5716

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2 (Leve Item ID 5062)
$ws.Range("H2").Value = 128.90909
$ws.Range("I2").Value = 126.6
$ws.Range("J2").Value = 133.85715
$ws.Range("K2").Value = 126.6
$ws.Range("L2").Value = 133.85715
$ws.Range("M2").Value = -13.59999999999999
$ws.Range("N2").Value = -359.85715
# Row 46 (Leve Item ID 2078)
$ws.Range("H46").Value = 5624.75
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 2152.75
$ws.Range("J102").Value = 1200
$ws.Range("L102").Value = 1200
$ws.Range("N102").Value = -4444
# Row 107 (Leve Item ID 27802)
$ws.Range("H107").Value = 37038050
$ws.Range("I107").Value = 144.75
$ws.Range("K107").Value = 144.75
$ws.Range("M107").Value = 1775.25
# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 750
$ws.Range("I113").Value = 750
$ws.Range("K113").Value = 750
$ws.Range("M113").Value = 1420
# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 4398.7617
$ws.Range("I122").Value = 3183.7693
$ws.Range("J122").Value = 6373.125
$ws.Range("K122").Value = 9551.3079
$ws.Range("L122").Value = 19119.375
$ws.Range("M122").Value = -7101.3079
$ws.Range("N122").Value = -24019.375

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 43 (Leve Item ID 4314)
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 550
$ws.Range("J46").Value = 700
$ws.Range("L46").Value = 700
$ws.Range("N46").Value = -1076
# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 11497888
$ws.Range("I93").Value = 11908394
$ws.Range("J93").Value = 3700
$ws.Range("K93").Value = 11908394
$ws.Range("L93").Value = 3700
$ws.Range("M93").Value = -11907146
$ws.Range("N93").Value = -6196

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 491.5
$ws.Range("I107").Value = 333
$ws.Range("K107").Value = 999
$ws.Range("M107").Value = 921
# Row 114 (Leve Item ID 25978)
$ws.Range("H114").Value = 42233
$ws.Range("J114").Value = 42233
$ws.Range("L114").Value = 42233
$ws.Range("N114").Value = -50911
